# Add files via upload
# The "Stats" column (H) values were re-typed from loose, unquoted
# pseudo-lists to valid JSON-style lists (string items wrapped in double
# quotes), and the active cell selection moved from D12 to H9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = ' [["Strenght", 10, 100], ["Age", 1, 99]]'
$ws.Range("H3").Value = ' [["Strenght", 10, 100]]'
$ws.Range("H4").Value = ' ["Strenght", 10, 100]'

$ws.Range("H9").Select() | Out-Null
